# Auto-generated script to update Goblin_Profits market-price derived columns (H-N)
# across all 8 sheets, per the scheduled runner diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ALC_changes = @(
    @{Cell="H8"; Value=810.5294},
    @{Cell="I8"; Value=906.0769},
    @{Cell="K8"; Value=2718.2307},
    @{Cell="M8"; Value=-2579.2307},
    @{Cell="H28"; Value=6489.4375},
    @{Cell="I28"; Value=7369.0713},
    @{Cell="J28"; Value=332},
    @{Cell="K28"; Value=7369.0713},
    @{Cell="L28"; Value=332},
    @{Cell="M28"; Value=-6884.0713},
    @{Cell="N28"; Value=-1302},
    @{Cell="H69"; Value=333349000},
    @{Cell="I69"; Value=0},
    @{Cell="J69"; Value=333349000},
    @{Cell="K69"; Value=0},
    @{Cell="L69"; Value=1000047000},
    @{Cell="M69"; Value=$null},
    @{Cell="N69"; Value=-1000048748},
    @{Cell="H72"; Value=333349000},
    @{Cell="I72"; Value=0},
    @{Cell="J72"; Value=333349000},
    @{Cell="K72"; Value=0},
    @{Cell="L72"; Value=3000141000},
    @{Cell="M72"; Value=$null},
    @{Cell="N72"; Value=-3000149736},
    @{Cell="H100"; Value=4314.8096},
    @{Cell="I100"; Value=1893.2307},
    @{Cell="K100"; Value=1893.2307},
    @{Cell="M100"; Value=-1352.2307},
    @{Cell="H106"; Value=3678.1482},
    @{Cell="I106"; Value=3395.5652},
    @{Cell="J106"; Value=5303},
    @{Cell="K106"; Value=3395.5652},
    @{Cell="L106"; Value=5303},
    @{Cell="M106"; Value=-2764.5652},
    @{Cell="N106"; Value=-6565},
    @{Cell="H113"; Value=4289.0454},
    @{Cell="I113"; Value=3844.9412},
    @{Cell="J113"; Value=5799},
    @{Cell="K113"; Value=3844.9412},
    @{Cell="L113"; Value=5799},
    @{Cell="M113"; Value=-590.9412000000002},
    @{Cell="N113"; Value=-12307},
    @{Cell="H127"; Value=1416.5714},
    @{Cell="I127"; Value=1416.5714},
    @{Cell="K127"; Value=4249.7142},
    @{Cell="M127"; Value=710.2857999999997},
    @{Cell="H129"; Value=1526.7646},
    @{Cell="I129"; Value=1045},
    @{Cell="K129"; Value=3135},
    @{Cell="M129"; Value=1865},
    @{Cell="H138"; Value=1943.9814},
    @{Cell="J138"; Value=2046.25},
    @{Cell="L138"; Value=6138.75},
    @{Cell="N138"; Value=-16418.75}
)
foreach ($chg in $ALC_changes) {
    if ($null -eq $chg.Value) {
        $ws.Range($chg.Cell).ClearContents()
    } else {
        $ws.Range($chg.Cell).Value = $chg.Value
    }
}

$ws = $wb.Worksheets.Item("ARM")
$ARM_changes = @(
    @{Cell="H32"; Value=4816.1465},
    @{Cell="I32"; Value=4193.838},
    @{Cell="J32"; Value=10572.5},
    @{Cell="K32"; Value=4193.838},
    @{Cell="L32"; Value=10572.5},
    @{Cell="M32"; Value=-3906.838},
    @{Cell="N32"; Value=-11146.5},
    @{Cell="H131"; Value=94999},
    @{Cell="J131"; Value=94999},
    @{Cell="L131"; Value=94999},
    @{Cell="N131"; Value=-105079}
)
foreach ($chg in $ARM_changes) {
    if ($null -eq $chg.Value) {
        $ws.Range($chg.Cell).ClearContents()
    } else {
        $ws.Range($chg.Cell).Value = $chg.Value
    }
}

$ws = $wb.Worksheets.Item("BSM")
$BSM_changes = @(
    @{Cell="H105"; Value=3727},
    @{Cell="I105"; Value=3472.25},
    @{Cell="K105"; Value=3472.25},
    @{Cell="M105"; Value=-1725.25}
)
foreach ($chg in $BSM_changes) {
    if ($null -eq $chg.Value) {
        $ws.Range($chg.Cell).ClearContents()
    } else {
        $ws.Range($chg.Cell).Value = $chg.Value
    }
}

$ws = $wb.Worksheets.Item("CRP")
$CRP_changes = @(
    @{Cell="H10"; Value=1656.3334},
    @{Cell="I10"; Value=1488.375},
    @{Cell="J10"; Value=3000},
    @{Cell="K10"; Value=1488.375},
    @{Cell="L10"; Value=3000},
    @{Cell="M10"; Value=-1349.375},
    @{Cell="N10"; Value=-3278},
    @{Cell="H92"; Value=50601},
    @{Cell="J92"; Value=50601},
    @{Cell="L92"; Value=50601},
    @{Cell="N92"; Value=-55593},
    @{Cell="H132"; Value=2044.7059},
    @{Cell="I132"; Value=2046.9166},
    @{Cell="K132"; Value=6140.7498},
    @{Cell="M132"; Value=-3610.7498}
)
foreach ($chg in $CRP_changes) {
    if ($null -eq $chg.Value) {
        $ws.Range($chg.Cell).ClearContents()
    } else {
        $ws.Range($chg.Cell).Value = $chg.Value
    }
}

$ws = $wb.Worksheets.Item("CUL")
$CUL_changes = @(
    @{Cell="H2"; Value=148.6875},
    @{Cell="I2"; Value=142.125},
    @{Cell="J2"; Value=155.25},
    @{Cell="K2"; Value=852.75},
    @{Cell="L2"; Value=931.5},
    @{Cell="M2"; Value=-739.75},
    @{Cell="N2"; Value=-1157.5},
    @{Cell="H55"; Value=264700.9},
    @{Cell="J55"; Value=15474.632},
    @{Cell="L55"; Value=46423.896},
    @{Cell="N55"; Value=-46777.896},
    @{Cell="H122"; Value=1628.7646},
    @{Cell="I122"; Value=351.2},
    @{Cell="J122"; Value=2161.0833},
    @{Cell="K122"; Value=3160.8},
    @{Cell="L122"; Value=19449.7497},
    @{Cell="M122"; Value=-710.7999999999997},
    @{Cell="N122"; Value=-24349.7497},
    @{Cell="H130"; Value=3799.5386},
    @{Cell="I130"; Value=3427.7144},
    @{Cell="J130"; Value=4233.3335},
    @{Cell="K130"; Value=10283.1432},
    @{Cell="L130"; Value=12700.0005},
    @{Cell="M130"; Value=-5263.143199999999},
    @{Cell="N130"; Value=-22740.0005}
)
foreach ($chg in $CUL_changes) {
    if ($null -eq $chg.Value) {
        $ws.Range($chg.Cell).ClearContents()
    } else {
        $ws.Range($chg.Cell).Value = $chg.Value
    }
}

$ws = $wb.Worksheets.Item("GSM")
$GSM_changes = @(
    @{Cell="H9"; Value=1525},
    @{Cell="I9"; Value=412.5},
    @{Cell="J9"; Value=3750},
    @{Cell="K9"; Value=412.5},
    @{Cell="L9"; Value=3750},
    @{Cell="M9"; Value=-242.5},
    @{Cell="N9"; Value=-4090},
    @{Cell="H69"; Value=0},
    @{Cell="J69"; Value=0},
    @{Cell="L69"; Value=0},
    @{Cell="N69"; Value=$null},
    @{Cell="H72"; Value=0},
    @{Cell="J72"; Value=0},
    @{Cell="L72"; Value=0},
    @{Cell="N72"; Value=$null},
    @{Cell="H102"; Value=3772.8857},
    @{Cell="I102"; Value=2048.2593},
    @{Cell="J102"; Value=9593.5},
    @{Cell="K102"; Value=2048.2593},
    @{Cell="L102"; Value=9593.5},
    @{Cell="M102"; Value=-426.2593000000002},
    @{Cell="N102"; Value=-12837.5},
    @{Cell="H113"; Value=35721176},
    @{Cell="I113"; Value=76926536},
    @{Cell="J113"; Value=9866.467000000001},
    @{Cell="K113"; Value=76926536},
    @{Cell="L113"; Value=9866.467000000001},
    @{Cell="M113"; Value=-76924366},
    @{Cell="N113"; Value=-14206.467},
    @{Cell="H122"; Value=40666.168},
    @{Cell="I122"; Value=42799.4},
    @{Cell="K122"; Value=128398.2},
    @{Cell="M122"; Value=-125948.2}
)
foreach ($chg in $GSM_changes) {
    if ($null -eq $chg.Value) {
        $ws.Range($chg.Cell).ClearContents()
    } else {
        $ws.Range($chg.Cell).Value = $chg.Value
    }
}

$ws = $wb.Worksheets.Item("LTW")
$LTW_changes = @(
    @{Cell="H61"; Value=5112},
    @{Cell="I61"; Value=4277.077},
    @{Cell="J61"; Value=6468.75},
    @{Cell="K61"; Value=4277.077},
    @{Cell="L61"; Value=6468.75},
    @{Cell="M61"; Value=-4075.077},
    @{Cell="N61"; Value=-6872.75},
    @{Cell="H64"; Value=26049.666},
    @{Cell="J64"; Value=26049.666},
    @{Cell="L64"; Value=26049.666},
    @{Cell="N64"; Value=-26499.666},
    @{Cell="H67"; Value=26049.666},
    @{Cell="J67"; Value=26049.666},
    @{Cell="L67"; Value=26049.666},
    @{Cell="N67"; Value=-27609.666},
    @{Cell="H93"; Value=6159.9473},
    @{Cell="J93"; Value=6822.727},
    @{Cell="L93"; Value=6822.727},
    @{Cell="N93"; Value=-9318.726999999999},
    @{Cell="H113"; Value=5112},
    @{Cell="I113"; Value=4277.077},
    @{Cell="J113"; Value=6468.75},
    @{Cell="K113"; Value=4277.077},
    @{Cell="L113"; Value=6468.75},
    @{Cell="M113"; Value=-2107.077},
    @{Cell="N113"; Value=-10808.75}
)
foreach ($chg in $LTW_changes) {
    if ($null -eq $chg.Value) {
        $ws.Range($chg.Cell).ClearContents()
    } else {
        $ws.Range($chg.Cell).Value = $chg.Value
    }
}

$ws = $wb.Worksheets.Item("WVR")
$WVR_changes = @(
    @{Cell="H122"; Value=2606.625},
    @{Cell="I122"; Value=1328.2963},
    @{Cell="K122"; Value=3984.8889},
    @{Cell="M122"; Value=-1534.8889}
)
foreach ($chg in $WVR_changes) {
    if ($null -eq $chg.Value) {
        $ws.Range($chg.Cell).ClearContents()
    } else {
        $ws.Range($chg.Cell).Value = $chg.Value
    }
}

Write-Host "Applied 173 cell updates across 8 sheets."
